# fall 23 week 14 inputs
# Appends 33 new rows (1305-1337) of matchup data to the "Nine" sheet,
# mirroring the four columns A-D already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$startRow = 1305
$numRows = 33
$numCols = 4
$endRow = $startRow + $numRows - 1

$data = New-Object 'object[,]' $numRows,$numCols
$data[0,0] = 6
$data[0,1] = 16
$data[0,2] = 7
$data[0,3] = 4
$data[1,0] = 3
$data[1,1] = 12
$data[1,2] = 4
$data[1,3] = 8
$data[2,0] = 6
$data[2,1] = 13
$data[2,2] = 5
$data[2,3] = 7
$data[3,0] = 5
$data[3,1] = 7
$data[3,2] = 4
$data[3,3] = 13
$data[4,0] = 4
$data[4,1] = 12
$data[4,2] = 1
$data[4,3] = 8
$data[5,0] = 4
$data[5,1] = 6
$data[5,2] = 3
$data[5,3] = 14
$data[6,0] = 4
$data[6,1] = 3
$data[6,2] = 3
$data[6,3] = 17
$data[7,0] = 5
$data[7,1] = 7
$data[7,2] = 2
$data[7,3] = 13
$data[8,0] = 5
$data[8,1] = 12
$data[8,2] = 3
$data[8,3] = 8
$data[9,0] = 5
$data[9,1] = 2
$data[9,2] = 3
$data[9,3] = 18
$data[10,0] = 4
$data[10,1] = 17
$data[10,2] = 2
$data[10,3] = 3
$data[11,0] = 6
$data[11,1] = 13
$data[11,2] = 9
$data[11,3] = 7
$data[12,0] = 3
$data[12,1] = 14
$data[12,2] = 2
$data[12,3] = 6
$data[13,0] = 5
$data[13,1] = 3
$data[13,2] = 4
$data[13,3] = 17
$data[14,0] = 4
$data[14,1] = 4
$data[14,2] = 3
$data[14,3] = 16
$data[15,0] = 5
$data[15,1] = 8
$data[15,2] = 4
$data[15,3] = 12
$data[16,0] = 4
$data[16,1] = 14
$data[16,2] = 5
$data[16,3] = 6
$data[17,0] = 4
$data[17,1] = 13
$data[17,2] = 5
$data[17,3] = 7
$data[18,0] = 5
$data[18,1] = 2
$data[18,2] = 7
$data[18,3] = 18
$data[19,0] = 4
$data[19,1] = 8
$data[19,2] = 3
$data[19,3] = 12
$data[20,0] = 5
$data[20,1] = 7
$data[20,2] = 3
$data[20,3] = 13
$data[21,0] = 2
$data[21,1] = 7
$data[21,2] = 4
$data[21,3] = 13
$data[22,0] = 3
$data[22,1] = 16
$data[22,2] = 4
$data[22,3] = 4
$data[23,0] = 3
$data[23,1] = 8
$data[23,2] = 5
$data[23,3] = 12
$data[24,0] = 4
$data[24,1] = 4
$data[24,2] = 2
$data[24,3] = 16
$data[25,0] = 4
$data[25,1] = 12
$data[25,2] = 5
$data[25,3] = 8
$data[26,0] = 2
$data[26,1] = 13
$data[26,2] = 5
$data[26,3] = 7
$data[27,0] = 3
$data[27,1] = 13
$data[27,2] = 2
$data[27,3] = 7
$data[28,0] = 3
$data[28,1] = 13
$data[28,2] = 4
$data[28,3] = 7
$data[29,0] = 4
$data[29,1] = 15
$data[29,2] = 8
$data[29,3] = 5
$data[30,0] = 4
$data[30,1] = 14
$data[30,2] = 7
$data[30,3] = 6
$data[31,0] = 6
$data[31,1] = 7
$data[31,2] = 5
$data[31,3] = 13
$data[32,0] = 5
$data[32,1] = 8
$data[32,2] = 4
$data[32,3] = 12

$range = $ws.Range("A${startRow}:D${endRow}")
$range.Value = $data

$nextRow = $endRow + 1
$ws.Range("A$nextRow").Select()
